$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: merge the first two runs of the opening paragraph ("...bajo" + ". ")
# into a single run, while leaving "Los costos" / " directos" / " en que
# incurrira..." as their own separate runs (matching the target OOXML).
#
# This runtime recombines a paragraph's compatible runs whenever any edit
# touches it, so first trigger that recombination with a no-op Find/Replace
# over the "bajo. Los costos" span, then use a temporary bookmark (added and
# immediately removed) at each boundary we want to keep split -- that forces
# a clean run break at that exact character offset without leaving any
# bookmark behind.
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute("muy bajo. Los costos", $true, $false, $false, $false, $false, $true, 1, $false, "muy bajo. Los costos", 2)

$splitPoints = @(219, 229, 238)
foreach ($pos in $splitPoints) {
    $bm = $d.Range($pos, $pos)
    $d.Bookmarks.Add("TempSplit_$pos", $bm)
    $d.Bookmarks("TempSplit_$pos").Delete()
}

# ---------------------------------------------------------------------------
# Part 2: table - combine "$ " and "750" into a single run "$ 750" in the
# "Servidor Dedicado" row, then delete the whole "Conexion a Internet" row.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$cell = $t.Rows.Item(2).Cells.Item(2)
$null = $cell.Range.Find.Execute("$ 750", $true, $false, $false, $false, $false, $true, 1, $false, "$ 750", 2)

$t.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# Part 3: delete the "Para la contratacion del servicio de conexion a
# internet..." paragraph entirely (including its paragraph mark), and merge
# the runs of the following GSM paragraph back into one run (also drops the
# now-redundant _GoBack bookmark that split "por me" / "s.").
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Content.Paragraphs.Count; $i++) {
    $para = $d.Content.Paragraphs.Item($i)
    if ($para.Range.Text -like "Para la contratación*") {
        $target = $para
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

$null = $d.Content.Find.Execute("Para el envío de mensajes de texto se necesita contratar una línea de teléfono GSM, la cual incluya al menos 500 mensajes de texto por mes.", $true, $false, $false, $false, $false, $true, 1, $false, "Para el envío de mensajes de texto se necesita contratar una línea de teléfono GSM, la cual incluya al menos 500 mensajes de texto por mes.", 2)

# ---------------------------------------------------------------------------
# Part 4: "Rapidez ... " bullet - drop "un equipo con " before "acceso a
# internet", and move the _GoBack bookmark so it now splits the run right
# before " utilice el usuario ...".
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("un equipo con ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$rapidezPara = $null
for ($i = 1; $i -le $d.Content.Paragraphs.Count; $i++) {
    $para = $d.Content.Paragraphs.Item($i)
    if ($para.Range.Text -like "Rapidez*") {
        $rapidezPara = $para
        break
    }
}

$paraStart = $rapidezPara.Range.Start
$paraText = $rapidezPara.Range.Text

$colonIdx = $paraText.IndexOf(": al ser")
$splitPos = $paraStart + $colonIdx
$bm = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TempSplit_rapidez", $bm)
$d.Bookmarks("TempSplit_rapidez").Delete()

$bookmarkIdx = $paraText.IndexOf(" utilice el usuario")
$bookmarkPos = $paraStart + $bookmarkIdx
$bmGoBack = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmGoBack)

Write-Output "done"
